# Cyclic rotation of species-record data between rows 47, 48, 50 and 51
# (columns A, B, E, F, G, H, Q, R, AI). Mapping of the rotation, derived
# from the target diff:
#   new row 47 <- old row 48
#   new row 48 <- old row 50
#   new row 50 <- old row 51
#   new row 51 <- old row 47

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AI")

# Capture all current ("old") values first, since rows will be overwritten
# and we must not read already-modified data. Use Value2 (not Value) for
# reading, as Value has been observed to not reliably return scalar data.
$old = @{}
foreach ($r in 47, 48, 50, 51) {
    $old[$r] = @{}
    foreach ($col in $cols) {
        $old[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# Target row -> source row mapping.
$mapping = @{ 47 = 48; 48 = 50; 50 = 51; 51 = 47 }

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value = $old[$sourceRow][$col]
    }
}
